$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "24/10/2025"
$ws.Range("B7").Value = "Al Fateh"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Al-Ettifaq"
$ws.Range("F7").Value = "W"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 2.08
$ws.Range("L7").Value = 1.4
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 11
$ws.Range("O7").Value = 7
$ws.Range("P7").Value = 3
